# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across the
# ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 13960
$ws.Range("J21").Value = 13960
$ws.Range("L21").Value = 13960
$ws.Range("N21").Value = -14896
$ws.Range("H23").Value = 13960
$ws.Range("J23").Value = 13960
$ws.Range("L23").Value = 13960
$ws.Range("N23").Value = -14428
$ws.Range("H64").Value = 65155
$ws.Range("I64").Value = 127424.875
$ws.Range("J64").Value = 2885.125
$ws.Range("K64").Value = 127424.875
$ws.Range("L64").Value = 2885.125
$ws.Range("M64").Value = -127176.875
$ws.Range("N64").Value = -3381.125
$ws.Range("H67").Value = 65155
$ws.Range("I67").Value = 127424.875
$ws.Range("J67").Value = 2885.125
$ws.Range("K67").Value = 127424.875
$ws.Range("L67").Value = 2885.125
$ws.Range("M67").Value = -126566.875
$ws.Range("N67").Value = -4601.125
$ws.Range("H70").Value = 1139.3
$ws.Range("I70").Value = 1305.7142
$ws.Range("J70").Value = 751
$ws.Range("K70").Value = 3917.1426
$ws.Range("L70").Value = 2253
$ws.Range("M70").Value = -3647.1426
$ws.Range("N70").Value = -2793
$ws.Range("H73").Value = 1139.3
$ws.Range("I73").Value = 1305.7142
$ws.Range("J73").Value = 751
$ws.Range("K73").Value = 3917.1426
$ws.Range("L73").Value = 2253
$ws.Range("M73").Value = -2981.1426
$ws.Range("N73").Value = -4125
$ws.Range("H106").Value = 2693.8125
$ws.Range("I106").Value = 3033.2222
$ws.Range("J106").Value = 2257.4285
$ws.Range("K106").Value = 3033.2222
$ws.Range("L106").Value = 2257.4285
$ws.Range("M106").Value = -2402.2222
$ws.Range("N106").Value = -3519.4285
$ws.Range("H129").Value = 2782.75
$ws.Range("I129").Value = 6385.5293
$ws.Range("J129").Value = 1212.3077
$ws.Range("K129").Value = 19156.5879
$ws.Range("L129").Value = 3636.9231
$ws.Range("M129").Value = -14156.5879
$ws.Range("N129").Value = -13636.9231
$ws.Range("H132").Value = 6103136
$ws.Range("I132").Value = 6416066
$ws.Range("K132").Value = 19248198
$ws.Range("M132").Value = -19245668
$ws.Range("H137").Value = 1090.804
$ws.Range("I137").Value = 1119.8
$ws.Range("J137").Value = 985.36365
$ws.Range("K137").Value = 3359.4
$ws.Range("L137").Value = 2956.09095
$ws.Range("M137").Value = -809.3999999999996
$ws.Range("N137").Value = -8056.09095
$ws.Range("H141").Value = 1394.7377
$ws.Range("I141").Value = 1167.9286
$ws.Range("J141").Value = 3935
$ws.Range("K141").Value = 3503.7858
$ws.Range("L141").Value = 11805
$ws.Range("M141").Value = 1676.2142
$ws.Range("N141").Value = -22165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20443.076
$ws.Range("I32").Value = 4095.9875
$ws.Range("J32").Value = 129423.664
$ws.Range("K32").Value = 4095.9875
$ws.Range("L32").Value = 129423.664
$ws.Range("M32").Value = -3808.9875
$ws.Range("N32").Value = -129997.664
$ws.Range("H61").Value = 1656.5
$ws.Range("I61").Value = 987.10345
$ws.Range("J61").Value = 2678.2104
$ws.Range("K61").Value = 987.10345
$ws.Range("L61").Value = 2678.2104
$ws.Range("M61").Value = -775.10345
$ws.Range("N61").Value = -3102.2104
$ws.Range("H64").Value = 41250
$ws.Range("J64").Value = 41250
$ws.Range("L64").Value = 41250
$ws.Range("N64").Value = -41746
$ws.Range("H67").Value = 41250
$ws.Range("J67").Value = 41250
$ws.Range("L67").Value = 41250
$ws.Range("N67").Value = -42966
$ws.Range("H122").Value = 1216.909
$ws.Range("I122").Value = 1204.64
$ws.Range("J122").Value = 1255.25
$ws.Range("K122").Value = 3613.92
$ws.Range("L122").Value = 3765.75
$ws.Range("M122").Value = -1163.92
$ws.Range("N122").Value = -8665.75
$ws.Range("H136").Value = 1656.5
$ws.Range("I136").Value = 987.10345
$ws.Range("J136").Value = 2678.2104
$ws.Range("K136").Value = 2961.31035
$ws.Range("L136").Value = 8034.6312
$ws.Range("M136").Value = -411.3103499999997
$ws.Range("N136").Value = -13134.6312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 47996.668
$ws.Range("J62").Value = 47996.668
$ws.Range("L62").Value = 47996.668
$ws.Range("N62").Value = -49368.668
$ws.Range("H65").Value = 47996.668
$ws.Range("J65").Value = 47996.668
$ws.Range("L65").Value = 143990.004
$ws.Range("N65").Value = -150854.004
$ws.Range("H133").Value = 68000
$ws.Range("J133").Value = 68000
$ws.Range("L133").Value = 68000
$ws.Range("N133").Value = -78120
$ws.Range("H135").Value = 48125
$ws.Range("J135").Value = 48125
$ws.Range("L135").Value = 48125
$ws.Range("N135").Value = -58265
$ws.Range("H137").Value = 39996
$ws.Range("J137").Value = 39996
$ws.Range("L137").Value = 39996
$ws.Range("N137").Value = -50196

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 720.6
$ws.Range("J122").Value = 720.6
$ws.Range("L122").Value = 2161.8
$ws.Range("N122").Value = -7061.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1406.5
$ws.Range("I5").Value = 1054.2222
$ws.Range("J5").Value = 1573.3684
$ws.Range("K5").Value = 3162.6666
$ws.Range("L5").Value = 4720.1052
$ws.Range("M5").Value = -3050.6666
$ws.Range("N5").Value = -4944.1052
$ws.Range("H17").Value = 1500
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1500
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4500
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4838
$ws.Range("H122").Value = 494.16666
$ws.Range("I122").Value = 493
$ws.Range("K122").Value = 4437
$ws.Range("M122").Value = -1987
$ws.Range("H135").Value = 1406.5
$ws.Range("I135").Value = 1054.2222
$ws.Range("J135").Value = 1573.3684
$ws.Range("K135").Value = 9487.9998
$ws.Range("L135").Value = 14160.3156
$ws.Range("M135").Value = -6952.9998
$ws.Range("N135").Value = -19230.3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 23396.5
$ws.Range("I11").Value = 20006
$ws.Range("J11").Value = 24526.666
$ws.Range("K11").Value = 20006
$ws.Range("L11").Value = 24526.666
$ws.Range("M11").Value = -19866
$ws.Range("N11").Value = -24806.666
$ws.Range("H17").Value = 25404
$ws.Range("I17").Value = 1008
$ws.Range("J17").Value = 49800
$ws.Range("K17").Value = 1008
$ws.Range("L17").Value = 49800
$ws.Range("M17").Value = -838
$ws.Range("N17").Value = -50140
$ws.Range("H20").Value = 47006.332
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 47006.332
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 47006.332
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -47458.332
$ws.Range("H22").Value = 1234.4546
$ws.Range("I22").Value = 1249.6666
$ws.Range("J22").Value = 1231.0741
$ws.Range("K22").Value = 1249.6666
$ws.Range("L22").Value = 1231.0741
$ws.Range("M22").Value = -954.6666
$ws.Range("N22").Value = -1821.0741
$ws.Range("H27").Value = 1234.4546
$ws.Range("I27").Value = 1249.6666
$ws.Range("J27").Value = 1231.0741
$ws.Range("K27").Value = 1249.6666
$ws.Range("L27").Value = 1231.0741
$ws.Range("M27").Value = -1142.6666
$ws.Range("N27").Value = -1445.0741
$ws.Range("H68").Value = 6080.5
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 6996.6
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 6996.6
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -8494.6
$ws.Range("H71").Value = 6080.5
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 6996.6
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 34983
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -42471
$ws.Range("H132").Value = 2557.95
$ws.Range("I132").Value = 2660.3713
$ws.Range("K132").Value = 7981.113899999999
$ws.Range("M132").Value = -5451.113899999999
$ws.Range("H136").Value = 1245.2106
$ws.Range("I136").Value = 1106.4117
$ws.Range("K136").Value = 3319.2351
$ws.Range("M136").Value = -769.2351000000003

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 42602
$ws.Range("J14").Value = 42602
$ws.Range("L14").Value = 42602
$ws.Range("N14").Value = -42938
$ws.Range("H132").Value = 2446.195
$ws.Range("I132").Value = 2526.8918
$ws.Range("J132").Value = 1699.75
$ws.Range("K132").Value = 7580.6754
$ws.Range("L132").Value = 5099.25
$ws.Range("M132").Value = -5050.6754
$ws.Range("N132").Value = -10159.25
$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -55120
$ws.Range("H136").Value = 421.1091
$ws.Range("I136").Value = 267.86273
$ws.Range("J136").Value = 2375
$ws.Range("K136").Value = 803.5881899999999
$ws.Range("L136").Value = 7125
$ws.Range("M136").Value = 1746.41181
$ws.Range("N136").Value = -12225
